# Final approval of CRS Speech Recognition review items:
# mark the remaining "Open" review rows (H5, H12, H13) as "Closed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H5").Value = "Closed"
$ws.Range("H12").Value = "Closed"
$ws.Range("H13").Value = "Closed"

# Reflect the last-edited / selected cell as in the authored workbook.
$ws.Range("H13").Select()
